$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row's column captions: "<name>_old" -> "<name>_FV2404"
#    and "<name>_new" -> "<name>_FV2410" (the "diff" header in between is left
#    untouched).
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -like "*_old") {
        $cell.Value2 = $v -replace "_old$", "_FV2404"
    } elseif ($v -like "*_new") {
        $cell.Value2 = $v -replace "_new$", "_FV2410"
    }
}

# 2. Turn the used range into a real Excel Table ("Table1") so the renamed
#    headers double as column headers / filter buttons.
$lastRow = $ws.UsedRange.Rows.Count
$tableAddress = "A1:" + $ws.Cells.Item($lastRow, $lastCol).Address($false, $false)
$rng = $ws.Range($tableAddress)
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
